# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns.
# Source data are plain text cells (coinranking.com scrape), so values are
# written as literal text to avoid Excel's automatic numeric coercion
# (e.g. "0.7127" must stay text, not become the float 0.7127000000000001).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2';  Value = '29.353.60' },    @{ Cell = 'E2';  Value = '  +0.19%  ' },
    @{ Cell = 'D3';  Value = '1.878.94' },     @{ Cell = 'E3';  Value = '  +0.24%  ' },
                                                @{ Cell = 'E4';  Value = '  +0.02%  ' },
    @{ Cell = 'D5';  Value = '0.7127' },       @{ Cell = 'E5';  Value = '  +0.06%  ' },
    @{ Cell = 'D6';  Value = '242.25' },       @{ Cell = 'E6';  Value = '  -0.19%  ' },
                                                @{ Cell = 'E7';  Value = '  +0.06%  ' },
    @{ Cell = 'D8';  Value = '0.08061' },      @{ Cell = 'E8';  Value = '  +4.09%  ' },
    @{ Cell = 'D9';  Value = '0.3127' },       @{ Cell = 'E9';  Value = '  +0.54%  ' },
    @{ Cell = 'D10'; Value = '25.24' },        @{ Cell = 'E10'; Value = '  +0.55%  ' },
    @{ Cell = 'D11'; Value = '0.08353' },      @{ Cell = 'E11'; Value = '  -1.51%  ' },
    @{ Cell = 'D12'; Value = '1.873.21' },     @{ Cell = 'E12'; Value = '  +0.40%  ' },
    @{ Cell = 'D13'; Value = '5.254' },        @{ Cell = 'E13'; Value = '  +0.84%  ' },
    @{ Cell = 'D14'; Value = '0.7192' },       @{ Cell = 'E14'; Value = '  +1.11%  ' },
    @{ Cell = 'D15'; Value = '91.46' },        @{ Cell = 'E15'; Value = '  +0.08%  ' },
    @{ Cell = 'D16'; Value = '6.266' },        @{ Cell = 'E16'; Value = '  +4.76%  ' },
    @{ Cell = 'D17'; Value = '0.000008384' },  @{ Cell = 'E17'; Value = '  +0.94%  ' },
    @{ Cell = 'D18'; Value = '29.354.07' },    @{ Cell = 'E18'; Value = '  +0.19%  ' },
    @{ Cell = 'D19'; Value = '240.95' },       @{ Cell = 'E19'; Value = '  -0.69%  ' },
    @{ Cell = 'D20'; Value = '13.25' },        @{ Cell = 'E20'; Value = '  +0.27%  ' },
    @{ Cell = 'D21'; Value = '2.127.46' },     @{ Cell = 'E21'; Value = '  +0.05%  ' },
    @{ Cell = 'D22'; Value = '1.000' },        @{ Cell = 'E22'; Value = '  +0.10%  ' },
    @{ Cell = 'D23'; Value = '7.808' },        @{ Cell = 'E23'; Value = '  +0.12%  ' },
    @{ Cell = 'D25'; Value = '0.1587' },       @{ Cell = 'E25'; Value = '  -2.28%  ' },
    @{ Cell = 'D26'; Value = '163.26' },       @{ Cell = 'E26'; Value = '  +0.13%  ' },
    @{ Cell = 'D27'; Value = '9.070' },
    @{ Cell = 'D28'; Value = '18.55' },        @{ Cell = 'E28'; Value = '  +0.24%  ' },
    @{ Cell = 'D29'; Value = '1.507' },        @{ Cell = 'E29'; Value = '  -0.11%  ' },
    @{ Cell = 'D30'; Value = '4.422' },        @{ Cell = 'E30'; Value = '  +0.04%  ' },
    @{ Cell = 'D31'; Value = '4.338' },        @{ Cell = 'E31'; Value = '  +0.28%  ' },
    @{ Cell = 'D32'; Value = '1.201' },        @{ Cell = 'E32'; Value = '  -5.91%  ' },
    @{ Cell = 'D33'; Value = '0.05376' },      @{ Cell = 'E33'; Value = '  +2.22%  ' },
    @{ Cell = 'D34'; Value = '1.952' },        @{ Cell = 'E34'; Value = '  +1.62%  ' },
                                                @{ Cell = 'E35'; Value = '  +0.50%  ' },
    @{ Cell = 'D36'; Value = '0.7514' },       @{ Cell = 'E36'; Value = '  +0.87%  ' },
    @{ Cell = 'D37'; Value = '2.701' },        @{ Cell = 'E37'; Value = '  +0.63%  ' },
    @{ Cell = 'D38'; Value = '1.287.33' },     @{ Cell = 'E38'; Value = '  +10.69%  ' },
    @{ Cell = 'D39'; Value = '0.01883' },      @{ Cell = 'E39'; Value = '  +1.23%  ' },
                                                @{ Cell = 'E40'; Value = '  +0.80%  ' },
    @{ Cell = 'D41'; Value = '6.585' },        @{ Cell = 'E41'; Value = '  +3.68%  ' },
    @{ Cell = 'D42'; Value = '110.46' },       @{ Cell = 'E42'; Value = '  +3.34%  ' },
    @{ Cell = 'D43'; Value = '0.8921' },       @{ Cell = 'E43'; Value = '  +0.25%  ' },
    @{ Cell = 'D44'; Value = '73.19' },        @{ Cell = 'E44'; Value = '  +0.39%  ' },
                                                @{ Cell = 'E45'; Value = '  +8.86%  ' },
                                                @{ Cell = 'E46'; Value = '  +0.07%  ' },
    @{ Cell = 'D47'; Value = '2.018.96' },     @{ Cell = 'E47'; Value = '  -0.21%  ' },
                                                @{ Cell = 'E48'; Value = '  -0.12%  ' },
    @{ Cell = 'D49'; Value = '0.5203' },       @{ Cell = 'E49'; Value = '  +0.17%  ' },
    @{ Cell = 'D50'; Value = '9.472' },        @{ Cell = 'E50'; Value = '  +1.09%  ' },
                                                @{ Cell = 'E51'; Value = '  +1.53%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell[0] -eq 'D') {
        # Force text storage so values like "0.7127" / "25.24" aren't
        # auto-parsed into doubles (which would corrupt trailing digits).
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
